$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" row at the bottom of the table (row 29)
$ws.Range("A29").Value = "Total"

# Copy the formatting (bold, centered, bordered) from the row above so the
# new row matches the rest of the state-name column.
$ws.Range("A28").Copy()
$ws.Range("A29").PasteSpecial(-4122)

$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 36.22
$ws.Range("D29").Value = 10.41
$ws.Range("E29").Value = 81.89
$ws.Range("F29").Value = 44.45
$ws.Range("G29").Value = 99.45
$ws.Range("H29").Value = 86.78
$ws.Range("I29").Value = 74.67
$ws.Range("J29").Value = 40.23
$ws.Range("K29").Value = 10.82
$ws.Range("L29").Value = 10.17
$ws.Range("M29").Value = 51.31
$ws.Range("N29").Value = 8.28
$ws.Range("O29").Value = 70.86
$ws.Range("P29").Value = 99.38
